$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2 relabeling ---
$ws.Range("C2").Value = "Red (mailbox radio) "
$ws.Range("D2").Value = "Green (house radio)"
$ws.Range("E2").Value = "Yellow (house text)"

# --- Row 3: drop the "lights off & wait for RX Green (12)" note in E3 ---
$ws.Range("E3").Value = ""

# --- Row 8: drop stray "x" marker in E8 ---
$ws.Range("E8").Value = ""

# --- Row 13: introduce "dig pin HIGH" note in D13 and replace E13 ---
$ws.Range("D13").Value = "dig pin HIGH"
$ws.Range("E13").Value = "dig pin HIGH"

# --- Row height tweaks ---
$ws.Rows.Item(3).RowHeight = 23.85
$ws.Rows.Item(8).RowHeight = 23.85

# --- Update the saved selection/active cell ---
$ws.Range("H15").Select()
